$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New columns S (CurChapter) and T (CurStageNum) - Inventory MVP work adds
# two new tracked fields to the user DB table.
# ---------------------------------------------------------------------------

# Row 1 (field/header row) - copy formatting (bold+centered) from R1
$ws.Range("R1").Copy() | Out-Null
$ws.Range("S1:T1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("S1").Value = "CurChapter"
$ws.Range("T1").Value = "CurStageNum"

# Row 2 (type row) - copy formatting (vertical-center) from R2
$ws.Range("R2").Copy() | Out-Null
$ws.Range("S2:T2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("S2").Value = "int"
$ws.Range("T2").Value = "int"

# Row 3 (field name row, repeated) - copy formatting from R3
$ws.Range("R3").Copy() | Out-Null
$ws.Range("S3:T3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("S3").Value = "CurChapter"
$ws.Range("T3").Value = "CurStageNum"

# Row 4 (data row) - copy formatting from R4, then set values
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4:T4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 1

# ---------------------------------------------------------------------------
# Row 4 data edits for existing columns (rebalanced MVP stat defaults)
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 1
$ws.Range("O4").Value = "100"
$ws.Range("P4").Value = 100
$ws.Range("Q4").Value = 100

# ---------------------------------------------------------------------------
# Column widths (best-fit sizing applied by Excel when the new columns were
# added / table was widened)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.714286
$ws.Columns.Item(2).ColumnWidth = 10.285714
$ws.Columns.Item(3).ColumnWidth = 5.571429
$ws.Columns.Item(4).ColumnWidth = 6.857143
$ws.Columns.Item(5).ColumnWidth = 10.714286
$ws.Columns.Item(6).ColumnWidth = 4.142857
$ws.Columns.Item(7).ColumnWidth = 4.142857
$ws.Columns.Item(10).ColumnWidth = 14.857143
$ws.Columns.Item(11).ColumnWidth = 9.0
$ws.Columns.Item(12).ColumnWidth = 11.285714
$ws.Columns.Item(13).ColumnWidth = 9.714286
$ws.Columns.Item(18).ColumnWidth = 18.714286
$ws.Columns.Item(19).ColumnWidth = 11.285714
$ws.Columns.Item(20).ColumnWidth = 13.714286

# ---------------------------------------------------------------------------
# Selection moves to R7 (matches the author's last selection before saving)
# ---------------------------------------------------------------------------
$ws.Range("R7").Select() | Out-Null
